$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 0.3025064822817632
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.1296456352636128
$ws.Range("C4").Value = 1932
$ws.Range("D4").Value = 83.16831683168317
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 0.1721911321566939
$ws.Range("C6").Value = 2095
$ws.Range("D6").Value = 89.49167022639898
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = 0.7689021785561726
$ws.Range("C8").Value = 1936
$ws.Range("D8").Value = 83.1972496776966
$ws.Range("C9").Value = 2082
$ws.Range("D9").Value = 89.47142243231629
$ws.Range("C10").Value = 1594
$ws.Range("D10").Value = 68.32404629232748
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0.04286326618088298
$ws.Range("C12").Value = 1575
$ws.Range("D12").Value = 68.06395851339671
$ws.Range("C13").Value = 1176
$ws.Range("D13").Value = 50.82108902333622
$ws.Range("C15").Value = 2074
$ws.Range("D15").Value = 89.39655172413794
$ws.Range("C16").Value = 1242
$ws.Range("D16").Value = 53.14505776636713
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 0.7702182284980745
$ws.Range("C18").Value = 1197
$ws.Range("D18").Value = 51.24143835616438
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 0.1292546316242999
$ws.Range("C21").Value = 1930
$ws.Range("D21").Value = 83.15381301163292
$ws.Range("C23").Value = 2073
$ws.Range("D23").Value = 89.39197930142304
$ws.Range("C24").Value = 2072
$ws.Range("D24").Value = 89.38740293356342
$ws.Range("C25").Value = 1575
$ws.Range("D25").Value = 67.946505608283
$ws.Range("C27").Value = 1180
$ws.Range("D27").Value = 50.92792403970652
$ws.Range("C29").Value = 1310
$ws.Range("D29").Value = 55.72096980008507
$ws.Range("C30").Value = 2075
$ws.Range("D30").Value = 89.40112020680741
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 0.1723395088323998
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 0.1291989664082688
$ws.Range("C33").Value = 1931
$ws.Range("D33").Value = 83.16106804478898
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 0.1297016861219196
$ws.Range("C35").Value = 1576
$ws.Range("D35").Value = 68.13661910938175
$ws.Range("C37").Value = 1290
$ws.Range("D37").Value = 55.34105534105534
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = 0.1284796573875803
$ws.Range("C39").Value = 22
$ws.Range("D39").Value = 0.9421841541755889
$ws.Range("C40").Value = 1951
$ws.Range("D40").Value = 83.30486763450043
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 0.04269854824935952
$ws.Range("C42").Value = 1961
$ws.Range("D42").Value = 83.69611609048229
$ws.Range("C43").Value = 18
$ws.Range("D43").Value = 0.7682458386683738
$ws.Range("C44").Value = 1990
$ws.Range("D44").Value = 85.15190415062045
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 0.04278990158322636
$ws.Range("C46").Value = 3
$ws.Range("D46").Value = 0.1295336787564767
$ws.Range("C47").Value = 8
$ws.Range("D47").Value = 0.3454231433506045
$ws.Range("C49").Value = 1652
$ws.Range("D49").Value = 70.96219931271477
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = 0.1284796573875803
$ws.Range("C51").Value = 4
$ws.Range("D51").Value = 0.1713062098501071
$ws.Range("C52").Value = 1932
$ws.Range("D52").Value = 83.16831683168317
$ws.Range("C53").Value = 1182
$ws.Range("D53").Value = 50.88247955230306
$ws.Range("C55").Value = 1179
$ws.Range("D55").Value = 50.90673575129534
$ws.Range("C56").Value = 1934
$ws.Range("D56").Value = 83.36206896551724
$ws.Range("C57").Value = 1579
$ws.Range("D57").Value = 68.06034482758621
$ws.Range("C58").Value = 3
$ws.Range("D58").Value = 0.1284246575342466
$ws.Range("C59").Value = 21
$ws.Range("D59").Value = 0.898972602739726
$ws.Range("C61").Value = 5
$ws.Range("D61").Value = 0.2142245072836333
$ws.Range("C62").Value = 1615
$ws.Range("D62").Value = 69.19451585261353
$ws.Range("C63").Value = 18
$ws.Range("D63").Value = 0.7712082262210797
$ws.Range("C64").Value = 1242
$ws.Range("D64").Value = 53.14505776636713
$ws.Range("C65").Value = 18
$ws.Range("D65").Value = 0.7702182284980745
$ws.Range("C67").Value = 1290
$ws.Range("D67").Value = 55.34105534105534
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = 0.04282655246252676
$ws.Range("C69").Value = 2117
$ws.Range("D69").Value = 90.66381156316916
$ws.Range("C70").Value = 3
$ws.Range("D70").Value = 0.129757785467128
$ws.Range("C71").Value = 1573
$ws.Range("D71").Value = 68.0363321799308
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 0.04244482173174872
$ws.Range("C73").Value = 19
$ws.Range("D73").Value = 0.8064516129032258
$ws.Range("C74").Value = 2094
$ws.Range("D74").Value = 89.48717948717949
$ws.Range("C77").Value = 1309
$ws.Range("D77").Value = 55.70212765957447
$ws.Range("C78").Value = 9
$ws.Range("D78").Value = 0.38643194504079
$ws.Range("C80").Value = 2075
$ws.Range("D80").Value = 89.40112020680741
$ws.Range("C81").Value = 1179
$ws.Range("D81").Value = 50.79707022834985
$ws.Range("C83").Value = 1200
$ws.Range("D83").Value = 51.83585313174947
$ws.Range("C85").Value = 9
$ws.Range("D85").Value = 0.388768898488121
$ws.Range("C86").Value = 1575
$ws.Range("D86").Value = 68.06395851339671
$ws.Range("C87").Value = 4
$ws.Range("D87").Value = 0.1728608470181504
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 0.2568493150684931
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 0.04280821917808219
